# Apply updated TPM-derived NATMI metrics to the LR-pairs worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 58.62893533333334
    "H2" = 175.886806
    "I2" = 0.5702456571409142
    "J2" = 0.5702456571409142
    "M2" = 7.145781666666667
    "N2" = 21.437345
    "O2" = 0.1148763047483796
    "P2" = 0.1148763047483796
    "Q2" = 418.949571241119
    "R2" = 3770.54614117007
    "S2" = 0.06550771389115967
    "T2" = 0.06550771389115967
    "G3" = 58.62893533333334
    "H3" = 175.886806
    "I3" = 0.5702456571409142
    "J3" = 0.5702456571409142
    "O3" = 0.2979187698001963
    "P3" = 0.2979187698001963
    "Q3" = 1086.498570317517
    "R3" = 9778.487132857652
    "S3" = 0.1698868846593257
    "T3" = 0.1698868846593257
    "G4" = 58.62893533333334
    "H4" = 175.886806
    "I4" = 0.5702456571409142
    "J4" = 0.5702456571409142
    "M4" = 22.83059133333333
    "N4" = 68.49177399999999
    "O4" = 0.3670268824232265
    "P4" = 0.3670268824232265
    "Q4" = 1338.53326290376
    "R4" = 12046.79936613384
    "S4" = 0.2092954857558139
    "T4" = 0.2092954857558139
    "G5" = 58.62893533333334
    "H5" = 175.886806
    "I5" = 0.5702456571409142
    "J5" = 0.5702456571409142
    "M5" = 13.69598566666667
    "N5" = 41.087957
    "O5" = 0.2201780430281976
    "P5" = 0.2201780430281976
    "Q5" = 802.9810579772602
    "R5" = 7226.829521795342
    "S5" = 0.125555572834615
    "T5" = 0.125555572834615
    "I6" = 0.1389799721218762
    "J6" = 0.1389799721218763
    "M6" = 7.145781666666667
    "N6" = 21.437345
    "O6" = 0.1148763047483796
    "P6" = 0.1148763047483796
    "Q6" = 102.1061695120889
    "R6" = 918.9555256088
    "S6" = 0.01596550563139396
    "T6" = 0.01596550563139396
    "I7" = 0.1389799721218762
    "J7" = 0.1389799721218763
    "O7" = 0.2979187698001963
    "P7" = 0.2979187698001963
    "S7" = 0.04140474232141494
    "T7" = 0.04140474232141495
    "I8" = 0.1389799721218762
    "J8" = 0.1389799721218763
    "M8" = 22.83059133333333
    "N8" = 68.49177399999999
    "O8" = 0.3670268824232265
    "P8" = 0.3670268824232265
    "Q8" = 326.2266239698844
    "R8" = 2936.03961572896
    "S8" = 0.05100938588715917
    "T8" = 0.05100938588715918
    "I9" = 0.1389799721218762
    "J9" = 0.1389799721218763
    "M9" = 13.69598566666667
    "N9" = 41.087957
    "O9" = 0.2201780430281976
    "P9" = 0.2201780430281976
    "Q9" = 195.7021218041422
    "R9" = 1761.31909623728
    "S9" = 0.03060033828190817
    "T9" = 0.03060033828190818
    "G10" = 27.27518533333334
    "H10" = 81.82555600000001
    "I10" = 0.265288050953297
    "J10" = 0.2652880509532971
    "M10" = 7.145781666666667
    "N10" = 21.437345
    "O10" = 0.1148763047483796
    "P10" = 0.1148763047483796
    "Q10" = 194.9025193098689
    "R10" = 1754.12267378882
    "S10" = 0.03047531098741461
    "T10" = 0.03047531098741462
    "G11" = 27.27518533333334
    "H11" = 81.82555600000001
    "I11" = 0.265288050953297
    "J11" = 0.2652880509532971
    "O11" = 0.2979187698001963
    "P11" = 0.2979187698001963
    "Q11" = 505.4577522400169
    "R11" = 4549.119770160152
    "S11" = 0.07903428978269804
    "T11" = 0.07903428978269805
    "G12" = 27.27518533333334
    "H12" = 81.82555600000001
    "I12" = 0.265288050953297
    "J12" = 0.2652880509532971
    "M12" = 22.83059133333333
    "N12" = 68.49177399999999
    "O12" = 0.3670268824232265
    "P12" = 0.3670268824232265
    "Q12" = 622.7086098862604
    "R12" = 5604.377488976344
    "S12" = 0.09736784628552268
    "T12" = 0.0973678462855227
    "G13" = 27.27518533333334
    "H13" = 81.82555600000001
    "I13" = 0.265288050953297
    "J13" = 0.2652880509532971
    "M13" = 13.69598566666667
    "N13" = 41.087957
    "O13" = 0.2201780430281976
    "P13" = 0.2201780430281976
    "Q13" = 373.5605473810102
    "R13" = 3362.044926429092
    "S13" = 0.05841060389766171
    "T13" = 0.05841060389766173
    "E14" = 3
    "F14" = 1
    "G14" = 2.620337
    "H14" = 7.861011
    "I14" = 0.02548631978391236
    "J14" = 0.02548631978391236
    "M14" = 7.145781666666667
    "N14" = 21.437345
    "O14" = 0.1148763047483796
    "P14" = 0.1148763047483796
    "Q14" = 18.72435609508833
    "R14" = 168.519204855795
    "S14" = 0.002927774238411372
    "T14" = 0.002927774238411373
    "E15" = 3
    "F15" = 1
    "G15" = 2.620337
    "H15" = 7.861011
    "I15" = 0.02548631978391236
    "J15" = 0.02548631978391236
    "O15" = 0.2979187698001963
    "P15" = 0.2979187698001963
    "Q15" = 48.55951055675133
    "R15" = 437.035595010762
    "S15" = 0.007592853036757573
    "T15" = 0.007592853036757574
    "E16" = 3
    "F16" = 1
    "G16" = 2.620337
    "H16" = 7.861011
    "I16" = 0.02548631978391236
    "J16" = 0.02548631978391236
    "M16" = 22.83059133333333
    "N16" = 68.49177399999999
    "O16" = 0.3670268824232265
    "P16" = 0.3670268824232265
    "Q16" = 59.82384320261265
    "R16" = 538.4145888235139
    "S16" = 0.009354164494730751
    "T16" = 0.009354164494730753
    "E17" = 3
    "F17" = 1
    "G17" = 2.620337
    "H17" = 7.861011
    "I17" = 0.02548631978391236
    "J17" = 0.02548631978391236
    "M17" = 13.69598566666667
    "N17" = 41.087957
    "O17" = 0.2201780430281976
    "P17" = 0.2201780430281976
    "Q17" = 35.88809799383633
    "R17" = 322.992881944527
    "S17" = 0.005611528014012658
    "T17" = 0.005611528014012659
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Host "Updated $($updates.Count) cells with refreshed TPM values"
